$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column I header: meleeAtk -> physicalAtk (Style : Enemy 의 meleeAttack -> physicalAtk)
$ws.Range("I1").Value = "physicalAtk"

# Remove leftover test/scratch data outside the main table
$ws.Range("P1:P2").Clear()
$ws.Range("Q17:Q18").Clear()

# Highlight a subset of header cells (type, rank, element, size, attackAngle, rotationSpeed)
# with a bold, centered style on a new light accent fill.
$highlightCols = @("B1", "C1", "E1", "F1", "M1", "O1")
foreach ($addr in $highlightCols) {
    $rng = $ws.Range($addr)
    $rng.Font.Bold = $true
    $rng.Font.Size = 14
    $rng.Interior.Color = 15518374
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4108
}

# Update the last active selection/window view (cosmetic)
$ws.Range("U21").Select()
